$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "d9aa1b49-6548-445a-844d-efdbb3923b89"
$ws.Range("B9").Value = "Login with valid username and password"
$ws.Range("C9").Value = "PASSED"

$ws.Range("A10").Value = "e83e7079-1f20-4704-a1a2-abdd28c61a7d"
$ws.Range("B10").Value = "Create a new Country"
$ws.Range("C10").Value = "PASSED"

$ws.Range("A11").Value = "87a72036-24c3-4ee6-a83f-9bc671e194c5"
$ws.Range("B11").Value = "Delete a Country with parameters"
$ws.Range("C11").Value = "PASSED"
